# The commit simplifies the document's docDefaults (styles.xml):
#   - rPrDefault/rPr keeps only rFonts/sz/szCs/lang (drops the
#     redundant b/i/smallCaps/strike/color/u/shd/vertAlign entries,
#     all of which were just spelling out Word's own built-in
#     fallback values and had no visible effect).
#   - pPrDefault/pPr collapses down to a single
#     <w:spacing w:line="276" w:lineRule="auto"/> (drops keepNext/
#     keepLines/widowControl/pBdr/shd/spacing-before-after/ind/
#     contextualSpacing/jc, again all redundant copies of Word's
#     built-in defaults).
#
# The Word object model has no direct handle onto <w:docDefaults> -
# it is only reachable indirectly through the "Normal" style (every
# paragraph in this document already uses the Normal style, and
# Normal itself carries no direct formatting, so it inherits straight
# from docDefaults). Applying the same paragraph spacing through the
# Normal style's ParagraphFormat reproduces the exact effective
# formatting the diff leaves behind.

$d = $word.ActiveDocument
$normal = $d.Styles("Normal")

# pPrDefault collapses to exactly <w:spacing w:line="276" w:lineRule="auto"/>.
# wdLineSpaceMultiple = 5; LineSpacing is expressed in points and gets
# multiplied by 20 to land on twentieths-of-a-point ("276").
$normal.ParagraphFormat.LineSpacingRule = 5
$normal.ParagraphFormat.LineSpacing = 13.8
